# "Fruta / hortaliza, semanal" update:
# A new weekly price record for Locoto (Vega Modelo de Temuco) is inserted
# as row 34, pushing the existing rows 34-64 down to 35-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 34 (shifts rows 34:64 down to 35:65,
# carrying their formatting/styles with them - matches Excel's
# Rows.Insert default of xlShiftDown/xlFormatFromLeftOrAbove).
$ws.Rows(34).Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44966
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100112042
$ws.Range("G34").Value = "Locoto"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 90
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = 3000
$ws.Range("N34").Value = "$/kilo"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 3000
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
